# Add Test Case 5 (row 6) to the Guru99 live-project test-case sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the new test case row (row 6) ------------------------------
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Verify you can create account in E-commerce site and can share wishlist to other people using email"
$ws.Range("C6").Value = "1. Goto http://live.demoguru99.com`n2. Click on my account link`n3. Click Create account link and fill new user information exoect email id`n4. Click register`n5. verify registration is done`n6. goto Tv menu`n7. Add product in your wish list`n8. Click share wishlist`n9. In next page enter Email and a message and click share wishlist`n10. Check wishlist is shared"
$ws.Range("D6").Value = "product = LOG LCD"
$ws.Range("E6").Value = "1. Account registration done`n2. Wishlist Shared Successfully"

# Row 6 grows taller to fit the new (longer) content.
$ws.Rows.Item(6).RowHeight = 200

# --- Update the view state to match where the author left the cursor ----
$ws.Range("E7").Select()
